$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. consumptionAssets: new narrow "id" column (A), re-zoom view, new selection
# ---------------------------------------------------------------------------
$wsConsumption = $wb.Worksheets.Item(1)
$wsConsumption.Columns.Item(1).ColumnWidth = 6.285714285714286
$wsConsumption.Range("D17").Select()
$excel.ActiveWindow.Zoom = 115

# ---------------------------------------------------------------------------
# 2. storageAssets: add the new "Electric Heavy Goods Vehicle" storage asset
#    (row 14), re-zoom view, new selection
# ---------------------------------------------------------------------------
$wsStorage = $wb.Worksheets.Item(4)
$wsStorage.Range("A14").Value2 = 13
$wsStorage.Range("B14").Value2 = "EHGV"
$wsStorage.Range("C14").Value2 = "STORAGE"
$wsStorage.Range("D14").Value2 = "ELECTRIC_HEAVY_GOODS_VEHICLE"
$wsStorage.Range("E14").Value2 = 110
$wsStorage.Range("F14").Value2 = 0
$wsStorage.Range("G14").Value2 = 1
$wsStorage.Range("H14").Value2 = 0
$wsStorage.Range("I14").Value2 = 0
$wsStorage.Range("J14").Value2 = 0
$wsStorage.Range("K14").Value2 = 0
$wsStorage.Range("L14").Value2 = 500
$wsStorage.Range("M14").Value2 = 0
$wsStorage.Range("D10").Select()
$excel.ActiveWindow.Zoom = 115

# ---------------------------------------------------------------------------
# 3. conversionAssets: re-zoom view, new selection (no longer the active tab)
# ---------------------------------------------------------------------------
$wsConversion = $wb.Worksheets.Item(3)
$wsConversion.Range("D24").Select()
$excel.ActiveWindow.Zoom = 145

# ---------------------------------------------------------------------------
# 4. productionAssets: becomes the active tab, new selection
# ---------------------------------------------------------------------------
$wsProduction = $wb.Worksheets.Item(2)
$wsProduction.Range("E35").Select()
